# Updated cryptos list refresh (Price / Volume(1h) columns, plus a handful of
# row re-orderings where coins swapped rank) matching the scraped GitHub
# Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several "Price" cells hold plain decimal-looking text (e.g. "1.19", "34.00")
# that Excel would otherwise silently coerce to a Double on assignment,
# dropping formatting like trailing zeros. Force those cells to Text format
# first so the assigned strings are preserved verbatim, matching the source
# workbook's inline-string cells.
$textRefs = @("D6","D7","D8","D9","D12","D13","D15","D17","D19","D22","D23","D24","D25","D27","D28","D29","D31","D32","D34","D35","D36","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50")
foreach ($r in $textRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range('D2').Value = '95.592.24'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '3.623.04'
$ws.Range('E3').Value = '  -2.23%  '
$ws.Range('E4').Value = '  +25.30%  '
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = '222.93'
$ws.Range('E6').Value = '  -6.06%  '
$ws.Range('D7').Value = '642.23'
$ws.Range('E7').Value = '  -2.25%  '
$ws.Range('D8').Value = '0.423'
$ws.Range('E8').Value = '  -5.81%  '
$ws.Range('D9').Value = '1.19'
$ws.Range('E9').Value = '  +5.43%  '
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').Value = '3.621.96'
$ws.Range('E11').Value = '  -2.22%  '
$ws.Range('D12').Value = '52.02'
$ws.Range('E12').Value = '  +16.26%  '
$ws.Range('D13').Value = '0.217'
$ws.Range('E13').Value = '  +4.68%  '
$ws.Range('E14').Value = '  -7.79%  '
$ws.Range('D15').Value = '6.53'
$ws.Range('E15').Value = '  -4.51%  '
$ws.Range('D16').Value = '4.298.65'
$ws.Range('E16').Value = '  -2.19%  '
$ws.Range('D17').Value = '25.16'
$ws.Range('E17').Value = '  +33.98%  '
$ws.Range('D18').Value = '95.422.74'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').Value = '8.91'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +6.47%  '
$ws.Range('D21').Value = '3.620.43'
$ws.Range('E21').Value = '  -2.25%  '
$ws.Range('B22').Value = 'Hedera'
$ws.Range('C22').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D22').Value = '0.283'
$ws.Range('E22').Value = '  +32.51%  '
$ws.Range('B23').Value = 'Stellar'
$ws.Range('C23').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D23').Value = '0.535'
$ws.Range('E23').Value = '  -0.16%  '
$ws.Range('D24').Value = '138.45'
$ws.Range('E24').Value = '  +16.66%  '
$ws.Range('D25').Value = '535.11'
$ws.Range('E25').Value = '  +1.75%  '
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('D27').Value = '0.0000202'
$ws.Range('E27').Value = '  -9.73%  '
$ws.Range('D28').Value = '7.00'
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('D29').Value = '13.20'
$ws.Range('E29').Value = '  -2.08%  '
$ws.Range('D30').Value = '3.794.41'
$ws.Range('E30').Value = '  -2.81%  '
$ws.Range('D31').Value = '13.52'
$ws.Range('E31').Value = '  +6.54%  '
$ws.Range('D32').Value = '3.14'
$ws.Range('E32').Value = '  +3.79%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  +3.42%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '34.00'
$ws.Range('E35').Value = '  +3.27%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').Value = '0.639'
$ws.Range('E36').Value = '  +7.25%  '
$ws.Range('E37').Value = '  -2.60%  '
$ws.Range('E38').Value = '  +0.18%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.0560'
$ws.Range('E39').Value = '  +23.64%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '7.41'
$ws.Range('E40').Value = '  +8.75%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '8.56'
$ws.Range('E42').Value = '  -2.13%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '596.82'
$ws.Range('E43').Value = '  -6.46%  '
$ws.Range('E44').Value = '  +7.95%  '
$ws.Range('D45').Value = '0.504'
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('D46').Value = '41.30'
$ws.Range('E46').Value = '  +2.96%  '
$ws.Range('D47').Value = '0.161'
$ws.Range('E47').Value = '  -3.52%  '
$ws.Range('D48').Value = '2.01'
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('D49').Value = '9.44'
$ws.Range('E49').Value = '  +7.47%  '
$ws.Range('D50').Value = '233.04'
$ws.Range('E50').Value = '  +13.36%  '
$ws.Range('E51').Value = '  -2.09%  '
